# Added functionality to be able to check todays date against pay period date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: label cell
$ws.Range("A1").Value = "l"

# E1: pay period date
$ws.Range("E1").Value = 44711
$ws.Range("E1").NumberFormat = "mm-dd-yy"

# A4:A17: sequential run of dates to check against today's date
$dates = 44701,44702,44703,44704,44705,44706,44707,44708,44709,44710,44711,44712,44713,44714
for ($i = 0; $i -lt $dates.Length; $i++) {
  $row = 4 + $i
  $cell = $ws.Cells.Item($row, 1)
  $cell.Value = $dates[$i]
  $cell.NumberFormat = "mm-dd-yy"
}

# Fit the date columns to their new content
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null

# Leave the selection on E2, matching the last saved view
[void]$ws.Range("E2").Select()
